$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44799
$ws.Range("I2").Value = "Primera"
$ws.Range("J2").Value = 800
$ws.Range("K2").Value = 1000
$ws.Range("L2").Value = 1200
$ws.Range("M2").Value = 1100
$ws.Range("P2").Value = 1100

# Row 3
$ws.Range("D3").Value = 44201
$ws.Range("I3").Value = "Segunda"
$ws.Range("J3").Value = 500
$ws.Range("K3").Value = 800
$ws.Range("L3").Value = 900
$ws.Range("M3").Value = 850
$ws.Range("P3").Value = 850

# Row 4
$ws.Range("D4").Value = 44174
$ws.Range("J4").Value = 800
$ws.Range("K4").Value = 450
$ws.Range("L4").Value = 500
$ws.Range("M4").Value = 475
$ws.Range("P4").Value = 475

# Row 5
$ws.Range("D5").Value = 44174
$ws.Range("I5").Value = "Tercera"
$ws.Range("J5").Value = 1200
$ws.Range("K5").Value = 250
$ws.Range("L5").Value = 350
$ws.Range("M5").Value = 300
$ws.Range("P5").Value = 300

# Row 6
$ws.Range("D6").Value = 44229
$ws.Range("J6").Value = 760
$ws.Range("K6").Value = 550
$ws.Range("L6").Value = 600
$ws.Range("M6").Value = 575
$ws.Range("P6").Value = 575

# Row 7
$ws.Range("D7").Value = 44267
$ws.Range("I7").Value = "Tercera"
$ws.Range("J7").Value = 400
$ws.Range("K7").Value = 500
$ws.Range("L7").Value = 600
$ws.Range("M7").Value = 550
$ws.Range("P7").Value = 550

# Row 8
$ws.Range("D8").Value = 44658
$ws.Range("I8").Value = "Segunda"
$ws.Range("J8").Value = 1000
$ws.Range("K8").Value = 600
$ws.Range("L8").Value = 650
$ws.Range("M8").Value = 625
$ws.Range("P8").Value = 625

# Row 10
$ws.Range("D10").Value = 44544
$ws.Range("J10").Value = 1000
$ws.Range("K10").Value = 600
$ws.Range("L10").Value = 650
$ws.Range("M10").Value = 625
$ws.Range("P10").Value = 625

# Row 11
$ws.Range("D11").Value = 44474
$ws.Range("J11").Value = 200
$ws.Range("K11").Value = 600
$ws.Range("L11").Value = 700
$ws.Range("M11").Value = 650
$ws.Range("P11").Value = 650

# Row 12
$ws.Range("D12").Value = 44245
$ws.Range("J12").Value = 800
$ws.Range("K12").Value = 850
$ws.Range("L12").Value = 900
$ws.Range("M12").Value = 875
$ws.Range("P12").Value = 875

# Row 13
$ws.Range("D13").Value = 44245
$ws.Range("I13").Value = "Segunda"
$ws.Range("J13").Value = 1000
$ws.Range("K13").Value = 750
$ws.Range("L13").Value = 800
$ws.Range("M13").Value = 775
$ws.Range("P13").Value = 775

# Row 14
$ws.Range("D14").Value = 44278
$ws.Range("J14").Value = 700
$ws.Range("L14").Value = 700
$ws.Range("M14").Value = 650
$ws.Range("P14").Value = 650

# Row 15
$ws.Range("D15").Value = 44278
$ws.Range("I15").Value = "Tercera"
$ws.Range("J15").Value = 400
$ws.Range("K15").Value = 500
$ws.Range("L15").Value = 600
$ws.Range("M15").Value = 550
$ws.Range("P15").Value = 550

# Row 16
$ws.Range("D16").Value = 44210
$ws.Range("J16").Value = 900
$ws.Range("K16").Value = 600
$ws.Range("L16").Value = 700
$ws.Range("M16").Value = 650
$ws.Range("P16").Value = 650

# Row 17
$ws.Range("D17").Value = 44224
$ws.Range("J17").Value = 800
$ws.Range("K17").Value = 850
$ws.Range("L17").Value = 900
$ws.Range("M17").Value = 875
$ws.Range("P17").Value = 875

# Row 18
$ws.Range("D18").Value = 44253
$ws.Range("I18").Value = "Segunda"
$ws.Range("J18").Value = 1000
$ws.Range("K18").Value = 800
$ws.Range("L18").Value = 900
$ws.Range("M18").Value = 850
$ws.Range("P18").Value = 850

# Row 19
$ws.Range("D19").Value = 44253
$ws.Range("J19").Value = 800
$ws.Range("K19").Value = 600
$ws.Range("L19").Value = 700
$ws.Range("M19").Value = 650
$ws.Range("P19").Value = 650
